$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4909.4736
$ws.Range("I64").Value = 3722.3333
$ws.Range("K64").Value = 3722.3333
$ws.Range("M64").Value = -3474.3333
$ws.Range("H67").Value = 4909.4736
$ws.Range("I67").Value = 3722.3333
$ws.Range("K67").Value = 3722.3333
$ws.Range("M67").Value = -2864.3333
$ws.Range("H70").Value = 896635.25
$ws.Range("I70").Value = 1870.4
$ws.Range("J70").Value = 1109674.5
$ws.Range("K70").Value = 5611.200000000001
$ws.Range("L70").Value = 3329023.5
$ws.Range("M70").Value = -5341.200000000001
$ws.Range("N70").Value = -3329563.5
$ws.Range("H73").Value = 896635.25
$ws.Range("I73").Value = 1870.4
$ws.Range("J73").Value = 1109674.5
$ws.Range("K73").Value = 5611.200000000001
$ws.Range("L73").Value = 3329023.5
$ws.Range("M73").Value = -4675.200000000001
$ws.Range("N73").Value = -3330895.5
$ws.Range("H80").Value = 896.4761999999999
$ws.Range("I80").Value = 361.83334
$ws.Range("J80").Value = 1609.3334
$ws.Range("K80").Value = 1085.50002
$ws.Range("L80").Value = 4828.0002
$ws.Range("M80").Value = -87.50001999999995
$ws.Range("N80").Value = -6824.0002
$ws.Range("H83").Value = 896.4761999999999
$ws.Range("I83").Value = 361.83334
$ws.Range("J83").Value = 1609.3334
$ws.Range("K83").Value = 3256.50006
$ws.Range("L83").Value = 14484.0006
$ws.Range("M83").Value = 1735.49994
$ws.Range("N83").Value = -24468.0006
$ws.Range("H138").Value = 4516.5537
$ws.Range("J138").Value = 4569.351
$ws.Range("L138").Value = 13708.053
$ws.Range("N138").Value = -23988.053

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 247.41667
$ws.Range("I5").Value = 298.8889
$ws.Range("K5").Value = 298.8889
$ws.Range("M5").Value = -186.8889
$ws.Range("H32").Value = 20795.117
$ws.Range("I32").Value = 20315.47
$ws.Range("K32").Value = 20315.47
$ws.Range("M32").Value = -20028.47
$ws.Range("H45").Value = 1722.0526
$ws.Range("I45").Value = 1654.7333
$ws.Range("K45").Value = 1654.7333
$ws.Range("M45").Value = -1277.7333
$ws.Range("H102").Value = 31252038
$ws.Range("I102").Value = 2134.8333
$ws.Range("K102").Value = 2134.8333
$ws.Range("M102").Value = -512.8332999999998
$ws.Range("H132").Value = 5007323
$ws.Range("I132").Value = 6672197
$ws.Range("K132").Value = 20016591
$ws.Range("M132").Value = -20014061

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 247.41667
$ws.Range("I4").Value = 298.8889
$ws.Range("K4").Value = 298.8889
$ws.Range("M4").Value = -183.8889
$ws.Range("H22").Value = 2289
$ws.Range("I22").Value = 1550.3334
$ws.Range("J22").Value = 4505
$ws.Range("K22").Value = 1550.3334
$ws.Range("L22").Value = 4505
$ws.Range("M22").Value = -1377.3334
$ws.Range("N22").Value = -4851
$ws.Range("H64").Value = 927.36365
$ws.Range("I64").Value = 1899.5
$ws.Range("J64").Value = 711.3333
$ws.Range("K64").Value = 1899.5
$ws.Range("L64").Value = 711.3333
$ws.Range("M64").Value = -1674.5
$ws.Range("N64").Value = -1161.3333
$ws.Range("H67").Value = 927.36365
$ws.Range("I67").Value = 1899.5
$ws.Range("J67").Value = 711.3333
$ws.Range("K67").Value = 1899.5
$ws.Range("L67").Value = 711.3333
$ws.Range("M67").Value = -1119.5
$ws.Range("N67").Value = -2271.3333
$ws.Range("H82").Value = 46442.5
$ws.Range("J82").Value = 99932
$ws.Range("L82").Value = 99932
$ws.Range("N82").Value = -100698
$ws.Range("H85").Value = 46442.5
$ws.Range("J85").Value = 99932
$ws.Range("L85").Value = 99932
$ws.Range("N85").Value = -102584
$ws.Range("H105").Value = 76944530
$ws.Range("I105").Value = 100026490
$ws.Range("J105").Value = 4673.6665
$ws.Range("K105").Value = 100026490
$ws.Range("L105").Value = 4673.6665
$ws.Range("M105").Value = -100024743
$ws.Range("N105").Value = -8167.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 813.3889
$ws.Range("I22").Value = 549.3077
$ws.Range("K22").Value = 549.3077
$ws.Range("M22").Value = -199.3077
$ws.Range("H31").Value = 55559588
$ws.Range("I31").Value = 80002000
$ws.Range("J31").Value = 11912427
$ws.Range("K31").Value = 80002000
$ws.Range("L31").Value = 11912427
$ws.Range("M31").Value = -80001705
$ws.Range("N31").Value = -11913017
$ws.Range("H34").Value = 55559588
$ws.Range("I34").Value = 80002000
$ws.Range("J34").Value = 11912427
$ws.Range("K34").Value = 80002000
$ws.Range("L34").Value = 11912427
$ws.Range("M34").Value = -80001798
$ws.Range("N34").Value = -11912831
$ws.Range("H105").Value = 2998
$ws.Range("I105").Value = 2998
$ws.Range("K105").Value = 2998
$ws.Range("M105").Value = -1251

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1476.3462
$ws.Range("J34").Value = 2288.2856
$ws.Range("L34").Value = 6864.8568
$ws.Range("N34").Value = -7032.8568
$ws.Range("H37").Value = 88883.336
$ws.Range("J37").Value = 88883.336
$ws.Range("L37").Value = 266650.008
$ws.Range("N37").Value = -266874.008
$ws.Range("H38").Value = 58823560
$ws.Range("I38").Value = 166666720
$ws.Range("J38").Value = 20.90909
$ws.Range("K38").Value = 500000160
$ws.Range("L38").Value = 62.72727
$ws.Range("M38").Value = -499999813
$ws.Range("N38").Value = -756.72727
$ws.Range("H39").Value = 2227.2942
$ws.Range("J39").Value = 3808
$ws.Range("L39").Value = 11424
$ws.Range("N39").Value = -12012
$ws.Range("H44").Value = 101365
$ws.Range("I44").Value = 1516.7778
$ws.Range("K44").Value = 4550.3334
$ws.Range("M44").Value = -4152.3334
$ws.Range("H51").Value = 2771.75
$ws.Range("I51").Value = 543.5
$ws.Range("J51").Value = 5000
$ws.Range("K51").Value = 1630.5
$ws.Range("L51").Value = 15000
$ws.Range("M51").Value = -1170.5
$ws.Range("N51").Value = -15920
$ws.Range("H55").Value = 1627.5385
$ws.Range("I55").Value = 921.25
$ws.Range("J55").Value = 1941.4445
$ws.Range("K55").Value = 2763.75
$ws.Range("L55").Value = 5824.333500000001
$ws.Range("M55").Value = -2586.75
$ws.Range("N55").Value = -6178.333500000001
$ws.Range("H80").Value = 2901
$ws.Range("I80").Value = 1702
$ws.Range("K80").Value = 5106
$ws.Range("M80").Value = -4170
$ws.Range("H83").Value = 2901
$ws.Range("I83").Value = 1702
$ws.Range("K83").Value = 15318
$ws.Range("M83").Value = -10638

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 28113.084
$ws.Range("I2").Value = 381.67856
$ws.Range("K2").Value = 381.67856
$ws.Range("M2").Value = -268.67856
$ws.Range("H46").Value = 32709.334
$ws.Range("J46").Value = 78949.5
$ws.Range("L46").Value = 78949.5
$ws.Range("N46").Value = -79261.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 915.3158
$ws.Range("I16").Value = 915.3158
$ws.Range("K16").Value = 915.3158
$ws.Range("M16").Value = -745.3158
$ws.Range("H22").Value = 3623.543
$ws.Range("I22").Value = 2233.4666
$ws.Range("K22").Value = 2233.4666
$ws.Range("M22").Value = -1938.4666
$ws.Range("H27").Value = 3623.543
$ws.Range("I27").Value = 2233.4666
$ws.Range("K27").Value = 2233.4666
$ws.Range("M27").Value = -2126.4666
$ws.Range("H61").Value = 2114.25
$ws.Range("I61").Value = 2553.0715
$ws.Range("K61").Value = 2553.0715
$ws.Range("M61").Value = -2351.0715
$ws.Range("H74").Value = 37857.5
$ws.Range("I74").Value = 28799
$ws.Range("K74").Value = 28799
$ws.Range("M74").Value = -27801
$ws.Range("H77").Value = 37857.5
$ws.Range("I77").Value = 28799
$ws.Range("K77").Value = 86397
$ws.Range("M77").Value = -81405
$ws.Range("H113").Value = 2114.25
$ws.Range("I113").Value = 2553.0715
$ws.Range("K113").Value = 2553.0715
$ws.Range("M113").Value = -383.0715

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1491.5454
$ws.Range("I100").Value = 741.2857
$ws.Range("K100").Value = 1482.5714
$ws.Range("M100").Value = -941.5714
